$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28

# Columns A and D look like a date / a zero-padded number respectively, so
# force them to be stored as text (leading apostrophe = Excel's "treat as
# text" input, same as a user typing '2025-01-16 into the cell).
$ws.Cells.Item($row, 1).Value = "'2025-01-16"
$ws.Cells.Item($row, 2).Value = "18:25:18"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "'02"
$ws.Cells.Item($row, 5).Value = 126803
$ws.Cells.Item($row, 6).Value = 141654
$ws.Cells.Item($row, 7).Value = 169407
$ws.Cells.Item($row, 8).Value = 158114
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142959
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191969
$ws.Cells.Item($row, 14).Value = 115510
$ws.Cells.Item($row, 15).Value = 45316
$ws.Cells.Item($row, 16).Value = 28538
$ws.Cells.Item($row, 17).Value = 65723
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 49281
$ws.Cells.Item($row, 20).Value = -1
